$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B4" = 0.6648627278815172
    "C4" = 0.667
    "D4" = 0.677166889836075
    "E4" = 0.6719999999999999
    "F4" = 0.5240434850327348
    "G4" = 0.53
    "H4" = 0.5204935166877871
    "I4" = 0.5205
    "J4" = 0.6434232903758921
    "K4" = 0.675
    "L4" = 0.6209477564719348
    "M4" = 0.6295000000000001

    "B5" = 0.4008757564102146
    "C5" = 0.351
    "D5" = 0.5758951126392987
    "E5" = 0.5945
    "F5" = 0.6680529450303134
    "H5" = 0.5099724786095357
    "I5" = 0.5175000000000001
    "J5" = 0.4035272835243034
    "K5" = 0.4029999999999999
    "L5" = 0.5745105902810784
    "M5" = 0.5669999999999999

    "B6" = 0.7218838531824505
    "C6" = 0.713
    "D6" = 0.7571914438514239
    "E6" = 0.7375
    "F6" = 0.5054654098681046
    "G6" = 0.506
    "H6" = 0.5074582709594921
    "I6" = 0.5054999999999999
    "J6" = 0.6999827819334058
    "K6" = 0.6910000000000001
    "L6" = 0.7396124781829192
    "M6" = 0.7209999999999999
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
